$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '26.180.90'
$ws.Cells.Item(2, 5).Value = '  -2.21%  '
$ws.Cells.Item(3, 4).Value = '1.668.51'
$ws.Cells.Item(3, 5).Value = '  -1.86%  '
$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = '1.005'
$ws.Cells.Item(4, 5).Value = '  +0.05%  '
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '217.02'
$ws.Cells.Item(5, 5).Value = '  -1.04%  '
$ws.Cells.Item(6, 5).Value = '  +0.57%  '
$ws.Cells.Item(7, 5).Value = '  +0.10%  '
$ws.Cells.Item(8, 5).Value = '  +1.20%  '
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '0.06402'
$ws.Cells.Item(9, 5).Value = '  +4.17%  '
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '21.58'
$ws.Cells.Item(10, 5).Value = '  -1.46%  '
$ws.Cells.Item(11, 5).Value = '  +1.18%  '
$ws.Cells.Item(12, 4).Value = '1.675.04'
$ws.Cells.Item(12, 5).Value = '  -1.38%  '
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '4.520'
$ws.Cells.Item(13, 5).Value = '  +1.72%  '
$ws.Cells.Item(14, 5).Value = '  +1.13%  '
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '0.000008572'
$ws.Cells.Item(15, 5).Value = '  +4.57%  '
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '64.25'
$ws.Cells.Item(16, 5).Value = '  -2.10%  '
$ws.Cells.Item(17, 4).Value = '26.218.17'
$ws.Cells.Item(17, 5).Value = '  -2.30%  '
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '4.935'
$ws.Cells.Item(18, 5).Value = '  -1.94%  '
$ws.Cells.Item(19, 5).Value = '  +0.09%  '
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '10.83'
$ws.Cells.Item(20, 5).Value = '  +0.82%  '
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '190.32'
$ws.Cells.Item(21, 5).Value = '  +2.50%  '
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '6.210'
$ws.Cells.Item(22, 5).Value = '  -0.57%  '
$ws.Cells.Item(23, 5).Value = '  +0.19%  '
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '145.35'
$ws.Cells.Item(24, 5).Value = '  -0.08%  '
$ws.Cells.Item(25, 5).Value = '  -0.65%  '
$ws.Cells.Item(26, 5).Value = '  +3.18%  '
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '15.63'
$ws.Cells.Item(27, 5).Value = '  +2.00%  '
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '0.06349'
$ws.Cells.Item(28, 5).Value = '  +12.23%  '
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '1.297'
$ws.Cells.Item(29, 5).Value = '  -2.14%  '
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '1.319'
$ws.Cells.Item(30, 5).Value = '  -0.95%  '
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '3.530'
$ws.Cells.Item(31, 5).Value = '  +1.45%  '
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '3.519'
$ws.Cells.Item(32, 5).Value = '  +2.03%  '
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '1.638'
$ws.Cells.Item(33, 5).Value = '  -1.85%  '
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '1.013'
$ws.Cells.Item(34, 5).Value = '  +0.34%  '
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '0.6072'
$ws.Cells.Item(35, 5).Value = '  +2.54%  '
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '2.372'
$ws.Cells.Item(36, 5).Value = '  -1.54%  '
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '2.656'
$ws.Cells.Item(37, 5).Value = '  +0.77%  '
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '6.153'
$ws.Cells.Item(38, 5).Value = '  +4.24%  '
$ws.Cells.Item(39, 5).Value = '  +0.31%  '
$ws.Cells.Item(40, 4).Value = '1.082.81'
$ws.Cells.Item(40, 5).Value = '  +1.15%  '
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '0.8646'
$ws.Cells.Item(41, 5).Value = '  +1.16%  '
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '1.010'
$ws.Cells.Item(42, 5).Value = '  +0.68%  '
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '101.16'
$ws.Cells.Item(43, 5).Value = '  +2.32%  '
$ws.Cells.Item(44, 4).Value = '1.818.46'
$ws.Cells.Item(44, 5).Value = '  -2.13%  '
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '0.00000000113'
$ws.Cells.Item(45, 5).Value = '  +6.61%  '
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '56.20'
$ws.Cells.Item(46, 5).Value = '  -0.71%  '
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '1.007'
$ws.Cells.Item(47, 5).Value = '  -0.54%  '
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '8.096'
$ws.Cells.Item(48, 5).Value = '  -0.57%  '
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '0.05205'
$ws.Cells.Item(49, 5).Value = '  -0.09%  '
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '0.4295'
$ws.Cells.Item(50, 5).Value = '  -0.86%  '
$ws.Cells.Item(51, 2).Value = 'Aptos'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '5.900'
$ws.Cells.Item(51, 5).Value = '  +4.35%  '
